$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 139 (this pushes the existing row 139 "SCALEXYZ..."
# wait - existing SCOPY-adjacent rows down by one, preserving all other data)
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new command entry
$ws.Cells.Item(139, 1).Value = "SCOPY"
$ws.Cells.Item(139, 2).Value = "Extracts a C3D ground section line to a polyline"

# Update the view to reflect where the author left the selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 125
$ws.Range("B140").Select()
